$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-15 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-16 Monday", 2) | Out-Null
$d.Content.Find.Execute("43×90=3870", $true, $false, $false, $false, $false, $true, 1, $false, "18×22=396", 2) | Out-Null
$d.Content.Find.Execute("50×95=4750", $true, $false, $false, $false, $false, $true, 1, $false, "95×53=5035", 2) | Out-Null
$d.Content.Find.Execute("90×28=2520", $true, $false, $false, $false, $false, $true, 1, $false, "55×71=3905", 2) | Out-Null
$d.Content.Find.Execute("59×63=3717", $true, $false, $false, $false, $false, $true, 1, $false, "43×51=2193", 2) | Out-Null
$d.Content.Find.Execute("80×42=3360", $true, $false, $false, $false, $false, $true, 1, $false, "56×74=4144", 2) | Out-Null
$d.Content.Find.Execute("71×54=3834", $true, $false, $false, $false, $false, $true, 1, $false, "62×66=4092", 2) | Out-Null
$d.Content.Find.Execute("48×96=4608", $true, $false, $false, $false, $false, $true, 1, $false, "36×25=900", 2) | Out-Null
$d.Content.Find.Execute("37×41=1517", $true, $false, $false, $false, $false, $true, 1, $false, "54×25=1350", 2) | Out-Null
$d.Content.Find.Execute("31×87=2697", $true, $false, $false, $false, $false, $true, 1, $false, "90×52=4680", 2) | Out-Null
$d.Content.Find.Execute("56×54=3024", $true, $false, $false, $false, $false, $true, 1, $false, "95×53=5035", 2) | Out-Null
$d.Content.Find.Execute("64×31=1984", $true, $false, $false, $false, $false, $true, 1, $false, "82×92=7544", 2) | Out-Null
$d.Content.Find.Execute("24×56=1344", $true, $false, $false, $false, $false, $true, 1, $false, "90×99=8910", 2) | Out-Null
$d.Content.Find.Execute("89×76=6764", $true, $false, $false, $false, $false, $true, 1, $false, "27×78=2106", 2) | Out-Null
$d.Content.Find.Execute("84×29=2436", $true, $false, $false, $false, $false, $true, 1, $false, "53×86=4558", 2) | Out-Null
$d.Content.Find.Execute("22×36=792", $true, $false, $false, $false, $false, $true, 1, $false, "12×86=1032", 2) | Out-Null
$d.Content.Find.Execute("66×66=4356", $true, $false, $false, $false, $false, $true, 1, $false, "40×66=2640", 2) | Out-Null
$d.Content.Find.Execute("48×36=1728", $true, $false, $false, $false, $false, $true, 1, $false, "22×30=660", 2) | Out-Null
$d.Content.Find.Execute("28×35=980", $true, $false, $false, $false, $false, $true, 1, $false, "99×86=8514", 2) | Out-Null
$d.Content.Find.Execute("42×33=1386", $true, $false, $false, $false, $false, $true, 1, $false, "25×33=825", 2) | Out-Null
$d.Content.Find.Execute("95×38=3610", $true, $false, $false, $false, $false, $true, 1, $false, "66×33=2178", 2) | Out-Null
$d.Content.Find.Execute("19×86=1634", $true, $false, $false, $false, $false, $true, 1, $false, "96×18=1728", 2) | Out-Null
$d.Content.Find.Execute("18×35=630", $true, $false, $false, $false, $false, $true, 1, $false, "54×92=4968", 2) | Out-Null
$d.Content.Find.Execute("71×88=6248", $true, $false, $false, $false, $false, $true, 1, $false, "85×95=8075", 2) | Out-Null
$d.Content.Find.Execute("53×55=2915", $true, $false, $false, $false, $false, $true, 1, $false, "53×66=3498", 2) | Out-Null
$d.Content.Find.Execute("81×19=1539", $true, $false, $false, $false, $false, $true, 1, $false, "21×75=1575", 2) | Out-Null
